$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/10/2025  Through  2/16/2025"

# --- Numeric cell updates (value + number format) ---
$ws.Range("F16").Value = 4
$ws.Range("F16").NumberFormat = "#,##0"
$ws.Range("G16").Value = 6
$ws.Range("G16").NumberFormat = "#,##0"
$ws.Range("H16").Value = -33.333333333333
$ws.Range("H16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I16").Value = 6
$ws.Range("I16").NumberFormat = "#,##0"
$ws.Range("J16").Value = 8
$ws.Range("J16").NumberFormat = "#,##0"
$ws.Range("K16").Value = -25
$ws.Range("K16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L16").Value = -45.454545454545
$ws.Range("L16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C17").Value = 6
$ws.Range("C17").NumberFormat = "#,##0"
$ws.Range("D17").Value = 6
$ws.Range("D17").NumberFormat = "#,##0"
$ws.Range("E17").Value = 0
$ws.Range("E17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F17").Value = 17
$ws.Range("F17").NumberFormat = "#,##0"
$ws.Range("G17").Value = 17
$ws.Range("G17").NumberFormat = "#,##0"
$ws.Range("I17").Value = 26
$ws.Range("I17").NumberFormat = "#,##0"
$ws.Range("J17").Value = 21
$ws.Range("J17").NumberFormat = "#,##0"
$ws.Range("K17").Value = 23.809523809523
$ws.Range("K17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L17").Value = -23.529411764705
$ws.Range("L17").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C18").Value = 7
$ws.Range("C18").NumberFormat = "#,##0"
$ws.Range("F18").Value = 12
$ws.Range("F18").NumberFormat = "#,##0"
$ws.Range("G18").Value = 4
$ws.Range("G18").NumberFormat = "#,##0"
$ws.Range("H18").Value = 200
$ws.Range("H18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I18").Value = 16
$ws.Range("I18").NumberFormat = "#,##0"
$ws.Range("K18").Value = 100
$ws.Range("K18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L18").Value = 166.666666666667
$ws.Range("L18").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C19").Value = 4
$ws.Range("C19").NumberFormat = "#,##0"
$ws.Range("D19").Value = 10
$ws.Range("D19").NumberFormat = "#,##0"
$ws.Range("E19").Value = -60
$ws.Range("E19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F19").Value = 14
$ws.Range("F19").NumberFormat = "#,##0"
$ws.Range("G19").Value = 31
$ws.Range("G19").NumberFormat = "#,##0"
$ws.Range("H19").Value = -54.838709677419
$ws.Range("H19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I19").Value = 35
$ws.Range("I19").NumberFormat = "#,##0"
$ws.Range("J19").Value = 50
$ws.Range("J19").NumberFormat = "#,##0"
$ws.Range("K19").Value = -30
$ws.Range("K19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L19").Value = -7.894736842105
$ws.Range("L19").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F20").Value = 9
$ws.Range("F20").NumberFormat = "#,##0"
$ws.Range("G20").Value = 1
$ws.Range("G20").NumberFormat = "#,##0"
$ws.Range("H20").Value = 800
$ws.Range("H20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I20").Value = 16
$ws.Range("I20").NumberFormat = "#,##0"
$ws.Range("K20").Value = 433.333333333333
$ws.Range("K20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L20").Value = 45.454545454545
$ws.Range("L20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C21").Value = 20
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("D21").Value = 17
$ws.Range("D21").NumberFormat = "#,##0"
$ws.Range("E21").Value = 17.647058823529
$ws.Range("E21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("F21").Value = 58
$ws.Range("F21").NumberFormat = "#,##0"
$ws.Range("G21").Value = 59
$ws.Range("G21").NumberFormat = "#,##0"
$ws.Range("H21").Value = -1.694915254237
$ws.Range("H21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("I21").Value = 102
$ws.Range("I21").NumberFormat = "#,##0"
$ws.Range("J21").Value = 90
$ws.Range("J21").NumberFormat = "#,##0"
$ws.Range("K21").Value = 13.333333333333
$ws.Range("K21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("L21").Value = 0
$ws.Range("L21").NumberFormat = "#,##0.00;""-""#,##0.00"
$ws.Range("C23").Value = 1
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("F23").Value = 2
$ws.Range("F23").NumberFormat = "#,##0"
$ws.Range("I23").Value = 2
$ws.Range("I23").NumberFormat = "#,##0"
$ws.Range("K23").Value = 100
$ws.Range("K23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L23").Value = -33.333333333333
$ws.Range("L23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C24").Value = 35
$ws.Range("C24").NumberFormat = "#,##0"
$ws.Range("D24").Value = 32
$ws.Range("D24").NumberFormat = "#,##0"
$ws.Range("E24").Value = 9.375
$ws.Range("E24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F24").Value = 134
$ws.Range("F24").NumberFormat = "#,##0"
$ws.Range("G24").Value = 127
$ws.Range("G24").NumberFormat = "#,##0"
$ws.Range("H24").Value = 5.511811023622
$ws.Range("H24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I24").Value = 208
$ws.Range("I24").NumberFormat = "#,##0"
$ws.Range("J24").Value = 210
$ws.Range("J24").NumberFormat = "#,##0"
$ws.Range("K24").Value = -0.95238095238
$ws.Range("K24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L24").Value = 27.60736196319
$ws.Range("L24").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C25").Value = 26
$ws.Range("C25").NumberFormat = "#,##0"
$ws.Range("D25").Value = 30
$ws.Range("D25").NumberFormat = "#,##0"
$ws.Range("E25").Value = -13.333333333333
$ws.Range("E25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F25").Value = 92
$ws.Range("F25").NumberFormat = "#,##0"
$ws.Range("G25").Value = 92
$ws.Range("G25").NumberFormat = "#,##0"
$ws.Range("H25").Value = 0
$ws.Range("H25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I25").Value = 143
$ws.Range("I25").NumberFormat = "#,##0"
$ws.Range("J25").Value = 137
$ws.Range("J25").NumberFormat = "#,##0"
$ws.Range("K25").Value = 4.379562043795
$ws.Range("K25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L25").Value = 38.83495145631
$ws.Range("L25").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C26").Value = 6
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 13
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -53.846153846153
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F26").Value = 26
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("G26").Value = 41
$ws.Range("G26").NumberFormat = "#,##0"
$ws.Range("H26").Value = -36.585365853658
$ws.Range("H26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I26").Value = 50
$ws.Range("I26").NumberFormat = "#,##0"
$ws.Range("J26").Value = 66
$ws.Range("J26").NumberFormat = "#,##0"
$ws.Range("K26").Value = -24.242424242424
$ws.Range("K26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L26").Value = -26.470588235294
$ws.Range("L26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C28").Value = 3
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 200
$ws.Range("E28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F28").Value = 7
$ws.Range("F28").NumberFormat = "#,##0"
$ws.Range("H28").Value = 133.333333333333
$ws.Range("H28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I28").Value = 9
$ws.Range("I28").NumberFormat = "#,##0"
$ws.Range("J28").Value = 5
$ws.Range("J28").NumberFormat = "#,##0"
$ws.Range("K28").Value = 80
$ws.Range("K28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L28").Value = 12.5
$ws.Range("L28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L29").Value = -100
$ws.Range("L29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L30").Value = -100
$ws.Range("L30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J31").Value = 3
$ws.Range("J31").NumberFormat = "#,##0"

# --- Cells converting from number back to placeholder text (copy value+format from donor cells) ---
$ws.Range("C14").Copy()
$ws.Range("G23").PasteSpecial(-4104)
$ws.Range("D18").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H23").PasteSpecial(-4104)
$ws.Range("D18").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("G27").PasteSpecial(-4104)
$ws.Range("D18").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("H27").PasteSpecial(-4104)
$ws.Range("D18").Copy()
$ws.Range("H27").PasteSpecial(-4122)

$excel.CutCopyMode = $false

